# ABC_XYZ fixed with time filter at 10/10/2022
# Append 19 new tyre-sales rows (135-153) to the "Holidays 2019" sheet,
# which holds the Tyres table in columns E:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

# Each entry: Tyre Size, Model, Param, Sales value, Date_of_sales (serial), Contragent
$data = @(
    @("315/80R22.5", "BEL-158M", "камневыт, груз, сер",   113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-158M", "камневыт, груз, трп",   113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-278",  "груз, сер",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-278",  "груз, трп",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-268",  "груз, сер",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-268",  "груз, трп",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-398",  "груз, сер",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-326",  "груз, сер",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-326",  "груз, трп",             113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-498",  "156L, груз, сер",       113, 44671, "БНХ ПОЛЬСКА"),
    @("315/80R22.5", "BEL-518",  "груз, сер",             113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "ИД-304М",  "16, груз, сер",         113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "ИД-304М",  "18, груз, сер",         113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "ИД-304М",  "16, груз, трп",         113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "ИД-304М",  "18, груз, трп",         113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "БИ-368М",  "18, груз, сер",         113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "БИ-368М",  "18, груз, сер",         113, 44671, "БНХ ПОЛЬСКА"),
    @("12.00R20",    "БИ-368М",  "18, груз, трп",         113, 44671, "БНХ ПОЛЬСКА"),
    @("195/65R15",   "BEL-337",  "б/к, легк, сер",        113, 44671, "БНХ ПОЛЬСКА")
)

$row = 135
foreach ($item in $data) {
    $ws.Cells.Item($row, 5).Value  = $item[0]
    $ws.Cells.Item($row, 6).Value  = $item[1]
    $ws.Cells.Item($row, 7).Value  = $item[2]
    $ws.Cells.Item($row, 8).Value  = $item[3]
    $ws.Cells.Item($row, 9).Value  = $item[4]
    $ws.Cells.Item($row, 9).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 10).Value = $item[5]
    $row++
}
